$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 31250314
$ws.Range("I11").Value = 31250314
$ws.Range("K11").Value = 31250314
$ws.Range("M11").Value = -31250174
$ws.Range("H17").Value = 4719741.5
$ws.Range("J17").Value = 4719741.5
$ws.Range("L17").Value = 14159224.5
$ws.Range("N17").Value = -14159560.5
$ws.Range("H112").Value = 19048994
$ws.Range("J112").Value = 22858662
$ws.Range("L112").Value = 68575986
$ws.Range("N112").Value = -68578202
$ws.Range("H132").Value = 4221338.5
$ws.Range("I132").Value = 1545.6857
$ws.Range("K132").Value = 4637.0571
$ws.Range("M132").Value = -2107.0571
$ws.Range("H138").Value = 3263.7764
$ws.Range("I138").Value = 1675.0968
$ws.Range("J138").Value = 4175.7964
$ws.Range("K138").Value = 5025.2904
$ws.Range("L138").Value = 12527.3892
$ws.Range("M138").Value = 114.7096000000001
$ws.Range("N138").Value = -22807.3892

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15510.4
$ws.Range("I32").Value = 11081.296
$ws.Range("J32").Value = 26354.068
$ws.Range("K32").Value = 11081.296
$ws.Range("L32").Value = 26354.068
$ws.Range("M32").Value = -10794.296
$ws.Range("N32").Value = -26928.068
$ws.Range("H45").Value = 151091.14
$ws.Range("I45").Value = 191616
$ws.Range("K45").Value = 191616
$ws.Range("M45").Value = -191239
$ws.Range("H74").Value = 14707514
$ws.Range("I74").Value = 1373.0555
$ws.Range("J74").Value = 31251922
$ws.Range("K74").Value = 1373.0555
$ws.Range("L74").Value = 31251922
$ws.Range("M74").Value = -499.0554999999999
$ws.Range("N74").Value = -31253670
$ws.Range("H77").Value = 14707514
$ws.Range("I77").Value = 1373.0555
$ws.Range("J77").Value = 31251922
$ws.Range("K77").Value = 6865.2775
$ws.Range("L77").Value = 156259610
$ws.Range("M77").Value = -2497.2775
$ws.Range("N77").Value = -156268346
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H121").Value = 32255
$ws.Range("J121").Value = 32255
$ws.Range("L121").Value = 32255
$ws.Range("N121").Value = -35749
$ws.Range("H132").Value = 8509.444
$ws.Range("I132").Value = 17000
$ws.Range("J132").Value = 7448.125
$ws.Range("K132").Value = 51000
$ws.Range("L132").Value = 22344.375
$ws.Range("M132").Value = -48470
$ws.Range("N132").Value = -27404.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H105").Value = 3946
$ws.Range("I105").Value = 2022.25
$ws.Range("J105").Value = 5869.75
$ws.Range("K105").Value = 2022.25
$ws.Range("L105").Value = 5869.75
$ws.Range("M105").Value = -275.25
$ws.Range("N105").Value = -9363.75
$ws.Range("H134").Value = 61999.176
$ws.Range("I134").Value = 3426.5715
$ws.Range("J134").Value = 335338
$ws.Range("K134").Value = 10279.7145
$ws.Range("L134").Value = 1006014
$ws.Range("M134").Value = -7744.7145
$ws.Range("N134").Value = -1011084

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 31252624
$ws.Range("J99").Value = 41669170
$ws.Range("L99").Value = 41669170
$ws.Range("N99").Value = -41672166
$ws.Range("H126").Value = 31252624
$ws.Range("J126").Value = 41669170
$ws.Range("L126").Value = 125007510
$ws.Range("N126").Value = -125012450
$ws.Range("H141").Value = 424373
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 424373
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 424373
$ws.Range("M141").Value = $null
$ws.Range("N141").Value = -434733

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9720.4
$ws.Range("I3").Value = 16280.25
$ws.Range("J3").Value = 5347.1665
$ws.Range("K3").Value = 48840.75
$ws.Range("L3").Value = 16041.4995
$ws.Range("M3").Value = -48728.75
$ws.Range("N3").Value = -16265.4995
$ws.Range("H5").Value = 2165986.2
$ws.Range("I5").Value = 600.7406999999999
$ws.Range("J5").Value = 6063680.5
$ws.Range("K5").Value = 1802.2221
$ws.Range("L5").Value = 18191041.5
$ws.Range("M5").Value = -1690.2221
$ws.Range("N5").Value = -18191265.5
$ws.Range("H86").Value = 1005.9048
$ws.Range("I86").Value = 1090.2222
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 3270.6666
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -2084.6666
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 1005.9048
$ws.Range("I89").Value = 1090.2222
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 9811.9998
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = -3883.9998
$ws.Range("N89").Value = -16356
$ws.Range("H107").Value = 513.9091
$ws.Range("I107").Value = 351.95
$ws.Range("J107").Value = 763.0769
$ws.Range("K107").Value = 1055.85
$ws.Range("L107").Value = 2289.2307
$ws.Range("M107").Value = 864.1500000000001
$ws.Range("N107").Value = -6129.2307
$ws.Range("H113").Value = 2143356.8
$ws.Range("I113").Value = 16666998
$ws.Range("J113").Value = 400519.9
$ws.Range("K113").Value = 50000994
$ws.Range("L113").Value = 1201559.7
$ws.Range("M113").Value = -49998824
$ws.Range("N113").Value = -1205899.7
$ws.Range("H116").Value = 2114.2856
$ws.Range("I116").Value = 1960
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 5880
$ws.Range("L116").Value = 7500
$ws.Range("M116").Value = -2438
$ws.Range("N116").Value = -14384
$ws.Range("H122").Value = 7213.7646
$ws.Range("I122").Value = 1112.2
$ws.Range("J122").Value = 15930.286
$ws.Range("K122").Value = 10009.8
$ws.Range("L122").Value = 143372.574
$ws.Range("M122").Value = -7559.800000000001
$ws.Range("N122").Value = -148272.574
$ws.Range("H123").Value = 9800
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 9800
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 29400
$ws.Range("M123").Value = $null
$ws.Range("N123").Value = -34300
$ws.Range("H130").Value = 8093.3335
$ws.Range("J130").Value = 8093.3335
$ws.Range("L130").Value = 24280.0005
$ws.Range("N130").Value = -34320.00049999999
$ws.Range("H131").Value = 2223316.2
$ws.Range("I131").Value = 6667355.5
$ws.Range("J131").Value = 1296.8334
$ws.Range("K131").Value = 20002066.5
$ws.Range("L131").Value = 3890.5002
$ws.Range("M131").Value = -19997026.5
$ws.Range("N131").Value = -13970.5002
$ws.Range("H135").Value = 2165986.2
$ws.Range("I135").Value = 600.7406999999999
$ws.Range("J135").Value = 6063680.5
$ws.Range("K135").Value = 5406.6663
$ws.Range("L135").Value = 54573124.5
$ws.Range("M135").Value = -2871.6663
$ws.Range("N135").Value = -54578194.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 51502.5
$ws.Range("I80").Value = 100005
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 100005
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -99007
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 51502.5
$ws.Range("I83").Value = 100005
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 500025
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -495033
$ws.Range("N83").Value = -24984
$ws.Range("H102").Value = 2816.2
$ws.Range("I102").Value = 2265.5
$ws.Range("K102").Value = 2265.5
$ws.Range("M102").Value = -643.5
$ws.Range("H126").Value = 10186.782
$ws.Range("I126").Value = 12155.611
$ws.Range("K126").Value = 36466.833
$ws.Range("M126").Value = -33996.833
$ws.Range("H132").Value = 18521332
$ws.Range("I132").Value = 83334830
$ws.Range("J132").Value = 3189.1428
$ws.Range("K132").Value = 250004490
$ws.Range("L132").Value = 9567.428400000001
$ws.Range("M132").Value = -250001960
$ws.Range("N132").Value = -14627.4284

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2123
$ws.Range("J22").Value = 2285.88
$ws.Range("L22").Value = 2285.88
$ws.Range("N22").Value = -2875.88
$ws.Range("H27").Value = 2123
$ws.Range("J27").Value = 2285.88
$ws.Range("L27").Value = 2285.88
$ws.Range("N27").Value = -2499.88
$ws.Range("H61").Value = 3852.5
$ws.Range("J61").Value = 3805
$ws.Range("L61").Value = 3805
$ws.Range("N61").Value = -4209
$ws.Range("H93").Value = 720.75
$ws.Range("I93").Value = 720.75
$ws.Range("K93").Value = 720.75
$ws.Range("M93").Value = 527.25
$ws.Range("H113").Value = 3852.5
$ws.Range("J113").Value = 3805
$ws.Range("L113").Value = 3805
$ws.Range("N113").Value = -8145
$ws.Range("H132").Value = 4130.129
$ws.Range("I132").Value = 3964.2593
$ws.Range("K132").Value = 11892.7779
$ws.Range("M132").Value = -9362.777900000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 5313
$ws.Range("J45").Value = 5313
$ws.Range("L45").Value = 5313
$ws.Range("N45").Value = -6295
$ws.Range("H132").Value = 1950.931
$ws.Range("I132").Value = 1015.4
$ws.Range("J132").Value = 2443.3157
$ws.Range("K132").Value = 3046.2
$ws.Range("L132").Value = 7329.9471
$ws.Range("M132").Value = -516.1999999999998
$ws.Range("N132").Value = -12389.9471
